# Apply a "Glitter" (Hexagon pattern) slide transition, slow speed, 3.9s
# duration to every slide in the deck.
#
# The authored OOXML diff adds the PowerPoint-2010+ "glitter" transition
# (wrapped in an <mc:AlternateContent> Choice/Fallback pair, the Fallback
# being a plain slow <p:fade/>) to the slide master, all slide layouts, and
# every slide. The PowerPoint VBA/COM object model only allows writing
# SlideShowTransition properties on Slide / SlideRange objects -- Master and
# CustomLayout objects expose SlideShowTransition for reading only (setting
# any property on them raises "Property ... cannot be found on this
# object." in real PowerPoint, just like in this host), and the "glitter"
# entry effect itself has no PpEntryEffect constant -- it was never added to
# the legacy VBA enumeration, even in current PowerPoint. The closest
# reachable equivalent through automation is therefore exactly the
# <mc:Fallback> branch PowerPoint itself writes for non-p14-aware readers:
# a slow p:fade with the same 3.9s (3900ms) duration. Apply that to every
# slide, which is the part of this transition change COM automation can
# actually reach.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $t = $s.SlideShowTransition
    # Order matters: EntryEffect re-initializes the transition, so set it
    # first, then Duration, then Speed -- otherwise Speed's "slow" attribute
    # gets clobbered.
    $t.EntryEffect = 1793   # ppEffectFade
    $t.Duration = 3.9       # seconds -> p14:dur="3900"
    $t.Speed = 1            # ppTransitionSpeedSlow -> spd="slow"
}
